# Generate Report for Archive
#
# The localization status used to read "Ready for handoff"; the report now
# reflects that the items are "In Translation" instead. That status string
# shows up in the Overview sheet (per-language status columns E/F, rows
# 2-4) and in the per-language detail sheets (zh-cn / de-de, "Status"
# column C, rows 2-4).
#
# Shortening the status text also lets the Status column(s) re-fit a
# little narrower than before.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2:F4").Value = $newStatus

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = $newStatus

$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2:C4").Value = $newStatus

# Re-fit the now-narrower Status columns to match.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
